$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the reference row: the shaded "SEO / Ameliorer la description" row.
# Several rows in the table share category "SEO" and even a similar
# meta-description link, so match on the distinctive wording of its
# "problem identified" cell instead.
$refRowIndex = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $c2 = $t.Cell($i, 2).Range.Text
    if ($c2 -like "*liorer la description*") {
        $refRowIndex = $i
        break
    }
}

$refRow = $t.Rows.Item($refRowIndex)

# Insert a brand-new row right before it; Word seeds it with the same
# shading/formatting as the surrounding row so we only need to touch what
# actually differs.
$newRow = $t.Rows.Add($refRow)

# The new row is the "PERFORMANCES" row -> white background instead of grey.
for ($c = 1; $c -le 6; $c++) {
    $newRow.Cells.Item($c).Shading.BackgroundPatternColor = 16777215
}

# --- Cell 1: category label ---
$newRow.Cells.Item(1).Range.Text = "PERFORMANCES"

# --- Cell 2: problem identified ---
$newRow.Cells.Item(2).Range.Text = "Optimisation script JS et CSS"

# --- Cell 3: explanation ---
$newRow.Cells.Item(3).Range.Text = "Certains fichier JS et CSS sont lourd et prennent du temps de chargement pour les pages."

# --- Cell 4: best practice / what was done (3 paragraphs + trailing blank one) ---
$para1 = "L’optimisation s’est déroulée en 2 étapes."
$para2 = "J’ai d’abord minimisé les fichiers les plus lourd afin d’alléger leur taille."
$para3 = "J’ai ensuite procédé à un pré chargement dans le head HTML afin que les script n’empêche pas l’affichage de la page et ainsi gagner du temps de chargement."
$newRow.Cells.Item(4).Range.Text = $para1 + [char]13 + $para2 + [char]13 + $para3 + [char]13

# --- Cell 5 ("Action recommandée" / X column): now empty ---
$newRow.Cells.Item(5).Range.Text = ""

# --- Cell 6 (reference URL column): now empty ---
$newRow.Cells.Item(6).Range.Text = ""
